$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the survey text in rows 1-3
$ws.Range("A1").Value = "Encuesta"
$ws.Range("A2").Value = "<i>Seleccione una opción</i>"
$ws.Range("A3").Value = "<em>Seleccione todas las opciones que correspondan</em>"

# Remove the now-unused rows 4 and 5 entirely so the used range shrinks to A1:A3
$ws.Range("A4:A5").EntireRow.Delete()
